# Shorten topic names in the schedule sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E15").Value = "Dynamic plots"
$ws.Range("E17").Value = "Dynamic plots (continuation)"
$ws.Range("E18").Value = "Dashboards"
$ws.Range("E20").Value = "Dashboards (continuation)"
$ws.Range("E23").Value = "Intro to R-Shiny"
$ws.Range("E25").Value = "R-Shiny part-II"
$ws.Range("E27").Value = "R-Shiny part-III"

# Update the active cell selection to match the saved view state
$ws.Range("E24").Select()
